# Updates cryptos list with latest price/volume data (GitHub Actions sync)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.866.84"
$ws.Range("E2").Value = "  -5.60%  "
$ws.Range("D3").Value = "2.578.55"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.64"
$ws.Range("E5").Value = "  -1.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.73"
$ws.Range("E6").Value = "  -3.89%  "
$ws.Range("E7").Value = "  -3.40%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.564"
$ws.Range("E9").Value = "  -2.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.05"
$ws.Range("E10").Value = "  -6.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0819"
$ws.Range("E11").Value = "  -3.17%  "
$ws.Range("E12").Value = "  -4.08%  "
$ws.Range("D13").Value = "2.977.85"
$ws.Range("E13").Value = "  -1.02%  "
$ws.Range("E14").Value = "  +1.41%  "
$ws.Range("D15").Value = "2.580.88"
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.893"
$ws.Range("E16").Value = "  -2.99%  "
$ws.Range("E17").Value = "  -3.77%  "
$ws.Range("D18").Value = "43.872.68"
$ws.Range("E18").Value = "  -5.72%  "
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").Value = "0.0₃0986"
$ws.Range("E20").Value = "  -2.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.49"
$ws.Range("E21").Value = "  -3.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.69"
$ws.Range("E22").Value = "  +2.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "266.47"
$ws.Range("E23").Value = "  -3.37%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.94"
$ws.Range("E24").Value = "  -2.93%  "
$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.23"
$ws.Range("E25").Value = "  +2.84%  "
$ws.Range("E26").Value = "  +1.88%  "
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("E28").Value = "  -3.27%  "
$ws.Range("E29").Value = "  -3.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.20"
$ws.Range("E30").Value = "  -2.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.21"
$ws.Range("E31").Value = "  -2.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.59"
$ws.Range("E32").Value = "  -1.14%  "
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.80"
$ws.Range("E34").Value = "  -2.32%  "
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "152.47"
$ws.Range("E35").Value = "  +0.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0816"
$ws.Range("E37").Value = "  -3.97%  "
$ws.Range("E38").Value = "  -1.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "24.27"
$ws.Range("E39").Value = "  +4.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.88"
$ws.Range("E40").Value = "  +5.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.58"
$ws.Range("E41").Value = "  -2.22%  "
$ws.Range("E42").Value = "  -4.63%  "
$ws.Range("E43").Value = "  -4.96%  "
$ws.Range("D44").Value = "2.041.23"
$ws.Range("E44").Value = "  -4.55%  "
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.07"
$ws.Range("E46").Value = "  -5.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.16"
$ws.Range("E47").Value = "  -3.74%  "
$ws.Range("E48").Value = "  +5.30%  "
$ws.Range("D49").Value = "2.835.71"
$ws.Range("E49").Value = "  -0.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "105.64"
$ws.Range("E50").Value = "  -3.12%  "
$ws.Range("E51").Value = "  -4.07%  "
